$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = 9.262999999999995
$ws.Range("B6").Value = 6.2615
$ws.Range("B7").Value = 5.296100000000003
$ws.Range("E7").Value = 15.92590000000001
$ws.Range("B8").Value = 6.089600000000002
$ws.Range("E11").Value = 17.17559999999999
$ws.Range("E12").Value = 17.8607
$ws.Range("E15").Value = 16.2922
$ws.Range("B16").Value = 6.8979
$ws.Range("B20").Value = 9.554499999999988
$ws.Range("E20").Value = 15.84209999999999
$ws.Range("B21").Value = 9.21689999999999
$ws.Range("E21").Value = 16.77150000000001
$ws.Range("E22").Value = 16.70390000000001
$ws.Range("E23").Value = 16.10269999999998
$ws.Range("B28").Value = 6.059500000000002
$ws.Range("B29").Value = 5.3853
$ws.Range("E29").Value = 17.22570000000002
$ws.Range("B30").Value = 5.133800000000003
$ws.Range("B32").Value = 7.404299999999995
$ws.Range("E34").Value = 17.099
$ws.Range("B40").Value = 8.880399999999996
$ws.Range("E42").Value = 16.5255
$ws.Range("E43").Value = 17.33220000000001
$ws.Range("E44").Value = 16.549
$ws.Range("E45").Value = 16.51959999999999
$ws.Range("B46").Value = 5.660100000000001
$ws.Range("E46").Value = 16.4997
$ws.Range("E50").Value = 16.6405
$ws.Range("B51").Value = 5.542099999999998
$ws.Range("E51").Value = 17.13800000000002
$ws.Range("B52").Value = 5.6088
$ws.Range("B57").Value = 5.3746
$ws.Range("E57").Value = 16.72029999999999
$ws.Range("B59").Value = 4.793899999999999
$ws.Range("B62").Value = 5.721300000000005
$ws.Range("E65").Value = 17.50380000000001
$ws.Range("B66").Value = 5.689299999999998
$ws.Range("E66").Value = 17.13650000000001
$ws.Range("E67").Value = 17.12650000000001
$ws.Range("B73").Value = 8.468
$ws.Range("B74").Value = 8.684899999999997
$ws.Range("B77").Value = 9.290899999999999
$ws.Range("E79").Value = 18.19510000000001
$ws.Range("E84").Value = 16.85349999999999
$ws.Range("E87").Value = 16.16629999999999
$ws.Range("B92").Value = 4.743600000000001
$ws.Range("E92").Value = 18.82100000000001
$ws.Range("E97").Value = 16.6641
$ws.Range("B100").Value = 5.745799999999998
